$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("M2").Value = "Delivered"
$ws.Range("M3").Value = "Delivered"
$ws.Range("M4").Value = "Delivered"
$ws.Range("M5").Value = "Delivered"
